$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Re-types the "<label>: <numbers>" suffix of a paragraph as its own run
# (leaving the "<label>: " prefix as the original run), matching how the
# author split these runs by retyping just the number list. The label
# prefix ("0: ", "1: ", ...) is kept untouched/unselected, only the part
# after it is replaced - which is exactly why PowerPoint ends up emitting
# two <a:r> runs for the paragraph instead of one.
function Split-ParagraphSuffix($textRange, $paraIndex, $newSuffix) {
    $para = $textRange.Paragraphs($paraIndex)
    $fullText = $para.Text
    $sepIdx = $fullText.IndexOf(": ")
    $prefixLen = $sepIdx + 2
    $len = $para.Length
    $suffix = $para.Characters($prefixLen + 1, $len - $prefixLen)
    $suffix.Text = $newSuffix
}

# --- Shape "TextBox 3" (id=4): P->v permutation table ---
$shp1 = $s.Shapes.Item("TextBox 3")
$tr1 = $shp1.TextFrame.TextRange

Split-ParagraphSuffix $tr1 2 "0,3,1"
Split-ParagraphSuffix $tr1 3 "1,3,2"
Split-ParagraphSuffix $tr1 4 "0,2,3"
Split-ParagraphSuffix $tr1 5 "0,1,2"

# --- Shape "TextBox 39" (id=40): e->v permutation table ---
$shp2 = $s.Shapes.Item("TextBox 39")
$tr2 = $shp2.TextFrame.TextRange

Split-ParagraphSuffix $tr2 3 "1,2"
